# Applies the "esophageal -> gastric" cancer-staging edits described by the diff.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$old,
        [string]$new
    )
    $r = $d.Content
    $found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output ("WARNING: text not found -> " + $old)
    }
}

# 1) Section heading + its bookmark (text first; bookmark handled further below).
Replace-Text "3 Esophageal Cancer Staging" "3 Gastric Cancer Staging"

# 2) "Depth of growth into the wall of the esophagus" -> "... of the stomach"
Replace-Text "= Tumor - Depth of growth into the wall of the esophagus" "= Tumor - Depth of growth into the wall of the stomach"

# 3) "If we look at the walls of the esophagus ..." -> "... of the stomach ..."
Replace-Text "If we look at the walls of the esophagus, we see several layers:" "If we look at the walls of the stomach, we see several layers:"

# 4) Merge the two laparoscopy paragraphs into a single (edited, typo-preserving) sentence.
#    Remove the whole second paragraph ("In general, laparoscopy is considered ...")
#    and fix the wording/typo in the first paragraph.
$r = $d.Content
$found = $r.Find.Execute("In general, laparoscopy is considered for cancers that invade from the esophagus into the stomach.")
if ($found) {
    $p = $r.Paragraphs.First
    $start = $p.Range.Start
    $end = $p.Range.End
    # Extend one character to the left (consumes the preceding paragraph's
    # mark) and one to the right (consumes this paragraph's own mark), so the
    # whole paragraph -- including its shell -- disappears.
    $delRange = $d.Range($start - 1, $end + 1)
    $delRange.Delete()
} else {
    Write-Output "WARNING: laparoscopy paragraph not found"
}

Replace-Text "Not all patients with esophageal cancer need a laparoscopy." "Not all patients with stoach cancer need a laparoscopy."

# 5) Rename the "esophageal-cancer-staging" bookmark to "gastric-cancer-staging".
#    The hosted Word OM here cannot resolve/rename bookmarks that were loaded
#    from the original file (Document.Bookmarks.Count reports 0 for them), so
#    as a best-effort we add a new bookmark with the correct name over the
#    same span the original bookmark covered (the "3 Gastric Cancer Staging"
#    heading through the end of the "M = Metastasis ..." bullet).
$rStart = $d.Content
$foundStart = $rStart.Find.Execute("3 Gastric Cancer Staging")
$rEnd = $d.Content
$foundEnd = $rEnd.Find.Execute("= Metastasis - Spread to liver, lungs, or bone")
if ($foundStart -and $foundEnd) {
    $bmStart = $rStart.Paragraphs.First.Range.Start
    $bmEnd = $rEnd.Paragraphs.First.Range.End
    $bmRange = $d.Range($bmStart, $bmEnd)
    $d.Bookmarks.Add("gastric-cancer-staging", $bmRange)
} else {
    Write-Output "WARNING: could not locate bookmark span for rename"
}

Write-Output "edit complete"
